$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1470588235294118
$ws.Range("C2").Value = 0.6911764705882353
$ws.Range("P2").Value = 0.07352941176470588
$ws.Range("S2").Value = 0.08823529411764706
$ws.Range("C3").Value = 0.06060606060606061
$ws.Range("J3").Value = 0.0202020202020202
$ws.Range("P3").Value = 0.7878787878787878
$ws.Range("S3").Value = 0.1313131313131313
$ws.Range("P4").Value = 0.6296296296296297
$ws.Range("S4").Value = 0.3703703703703703
$ws.Range("B6").Value = 0.03184713375796178
$ws.Range("D6").Value = 0.006369426751592357
$ws.Range("F6").Value = 0.03821656050955414
$ws.Range("J6").Value = 0.2292993630573248
$ws.Range("O6").Value = 0.03821656050955414
$ws.Range("Q6").Value = 0.2101910828025478
$ws.Range("R6").Value = 0.06369426751592357
$ws.Range("S6").Value = 0.3821656050955414
$ws.Range("B7").Value = 0.05625
$ws.Range("D7").Value = 0.0125
$ws.Range("F7").Value = 0.06875000000000001
$ws.Range("J7").Value = 0.13125
$ws.Range("O7").Value = 0.03125
$ws.Range("Q7").Value = 0.2625
$ws.Range("R7").Value = 0.05
$ws.Range("S7").Value = 0.3875
$ws.Range("B8").Value = 0.0389908256880734
$ws.Range("D8").Value = 0.01605504587155963
$ws.Range("F8").Value = 0.02981651376146789
$ws.Range("J8").Value = 0.0963302752293578
$ws.Range("O8").Value = 0.01146788990825688
$ws.Range("Q8").Value = 0.2362385321100917
$ws.Range("R8").Value = 0.07110091743119266
$ws.Range("S8").Value = 0.5
$ws.Range("B9").Value = 0.03861003861003861
$ws.Range("D9").Value = 0.003861003861003861
$ws.Range("E9").Value = 0.003861003861003861
$ws.Range("F9").Value = 0.03861003861003861
$ws.Range("J9").Value = 0.1235521235521236
$ws.Range("O9").Value = 0.0193050193050193
$ws.Range("Q9").Value = 0.250965250965251
$ws.Range("R9").Value = 0.05791505791505792
$ws.Range("S9").Value = 0.4633204633204633
$ws.Range("B10").Value = 0.07183364839319471
$ws.Range("D10").Value = 0.01890359168241966
$ws.Range("E10").Value = 0.000945179584120983
$ws.Range("F10").Value = 0.06994328922495274
$ws.Range("J10").Value = 0.1257088846880907
$ws.Range("O10").Value = 0.0113421550094518
$ws.Range("Q10").Value = 0.2448015122873346
$ws.Range("R10").Value = 0.07655954631379962
$ws.Range("S10").Value = 0.3799621928166352
$ws.Range("G11").Value = 0.1407766990291262
$ws.Range("J11").Value = 0.0970873786407767
$ws.Range("K11").Value = 0.1699029126213592
$ws.Range("L11").Value = 0.5825242718446602
$ws.Range("S11").Value = 0.009708737864077669
$ws.Range("G12").Value = 0.8064516129032258
$ws.Range("J12").Value = 0.1451612903225807
$ws.Range("L12").Value = 0.03225806451612903
$ws.Range("S12").Value = 0.01612903225806452
$ws.Range("G13").Value = 0.8461538461538461
$ws.Range("J13").Value = 0.1282051282051282
$ws.Range("S13").Value = 0.02564102564102564
$ws.Range("H15").Value = 0.1846153846153846
$ws.Range("I15").Value = 0.1179487179487179
$ws.Range("J15").Value = 0.3179487179487179
$ws.Range("K15").Value = 0.07179487179487179
$ws.Range("M15").Value = 0.03076923076923077
$ws.Range("O15").Value = 0.04102564102564103
$ws.Range("S15").Value = 0.2358974358974359
$ws.Range("F16").Value = 0.009523809523809525
$ws.Range("H16").Value = 0.1904761904761905
$ws.Range("I16").Value = 0.08571428571428572
$ws.Range("J16").Value = 0.4761904761904762
$ws.Range("K16").Value = 0.06666666666666667
$ws.Range("M16").Value = 0.01904761904761905
$ws.Range("O16").Value = 0.05714285714285714
$ws.Range("S16").Value = 0.09523809523809523
$ws.Range("F17").Value = 0.009980039920159681
$ws.Range("H17").Value = 0.2055888223552894
$ws.Range("I17").Value = 0.1277445109780439
$ws.Range("J17").Value = 0.407185628742515
$ws.Range("K17").Value = 0.06986027944111776
$ws.Range("M17").Value = 0.01996007984031936
$ws.Range("O17").Value = 0.06786427145708583
$ws.Range("S17").Value = 0.09181636726546906
$ws.Range("H18").Value = 0.2482758620689655
$ws.Range("I18").Value = 0.1517241379310345
$ws.Range("J18").Value = 0.3241379310344827
$ws.Range("K18").Value = 0.03448275862068965
$ws.Range("M18").Value = 0.02068965517241379
$ws.Range("O18").Value = 0.09655172413793103
$ws.Range("S18").Value = 0.1241379310344828
$ws.Range("F19").Value = 0.01974865350089767
$ws.Range("H19").Value = 0.2190305206463196
$ws.Range("I19").Value = 0.1292639138240574
$ws.Range("J19").Value = 0.3581687612208259
$ws.Range("K19").Value = 0.09515260323159784
$ws.Range("M19").Value = 0.01615798922800718
$ws.Range("N19").Value = 0.003590664272890485
$ws.Range("O19").Value = 0.06283662477558348
$ws.Range("S19").Value = 0.09605026929982047